$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F (copy formatting from the adjacent header cell)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$timestamps = @(
    "2021-10-05 10:51:17.753863",
    "2021-10-05 10:51:17.753875",
    "2021-10-05 10:51:17.753879",
    "2021-10-05 10:51:17.753883",
    "2021-10-05 10:51:17.753886",
    "2021-10-05 10:51:17.753889",
    "2021-10-05 10:51:17.753892",
    "2021-10-05 10:51:17.753895",
    "2021-10-05 10:51:17.753899",
    "2021-10-05 10:51:17.753902",
    "2021-10-05 10:51:17.753905",
    "2021-10-05 10:51:17.753908",
    "2021-10-05 10:51:17.753911",
    "2021-10-05 10:51:17.753914",
    "2021-10-05 10:51:17.753917",
    "2021-10-05 10:51:17.753920",
    "2021-10-05 10:51:17.753923",
    "2021-10-05 10:51:17.753926",
    "2021-10-05 10:51:17.753929",
    "2021-10-05 10:51:17.753932",
    "2021-10-05 10:51:17.753935",
    "2021-10-05 10:51:17.753948",
    "2021-10-05 10:51:17.753952",
    "2021-10-05 10:51:17.753955",
    "2021-10-05 10:51:17.753958",
    "2021-10-05 10:51:17.753962",
    "2021-10-05 10:51:17.753965",
    "2021-10-05 10:51:17.753968",
    "2021-10-05 10:51:17.753971",
    "2021-10-05 10:51:17.753974",
    "2021-10-05 10:51:17.753977",
    "2021-10-05 10:51:17.753979",
    "2021-10-05 10:51:17.753983",
    "2021-10-05 10:51:17.754006",
    "2021-10-05 10:51:17.754020",
    "2021-10-05 10:51:17.754024",
    "2021-10-05 10:51:17.754028",
    "2021-10-05 10:51:17.754031",
    "2021-10-05 10:51:17.754034",
    "2021-10-05 10:51:17.754038",
    "2021-10-05 10:51:17.754042",
    "2021-10-05 10:51:17.754045",
    "2021-10-05 10:51:17.754048",
    "2021-10-05 10:51:17.754051",
    "2021-10-05 10:51:17.754054",
    "2021-10-05 10:51:17.754057",
    "2021-10-05 10:51:17.754061",
    "2021-10-05 10:51:17.754064",
    "2021-10-05 10:51:17.754067",
    "2021-10-05 10:51:17.754070",
    "2021-10-05 10:51:17.754073",
    "2021-10-05 10:51:17.754076",
    "2021-10-05 10:51:17.754080",
    "2021-10-05 10:51:17.754083",
    "2021-10-05 10:51:17.754086",
    "2021-10-05 10:51:17.754089",
    "2021-10-05 10:51:17.754092",
    "2021-10-05 10:51:17.754095",
    "2021-10-05 10:51:17.754098",
    "2021-10-05 10:51:17.754101",
    "2021-10-05 10:51:17.754104",
    "2021-10-05 10:51:17.754107",
    "2021-10-05 10:51:17.754110",
    "2021-10-05 10:51:17.754113",
    "2021-10-05 10:51:17.754118",
    "2021-10-05 10:51:17.754122"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
